$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values look numeric ("7.00", "67.589.52", ...) but are
# stored as text in the source data. Force text format on exactly the cells
# being rewritten so Excel does not reinterpret them as numbers and drop
# significant trailing zeros / treat multi-dot strings as invalid numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '67.589.52'
$ws.Range("E2").Value = '  +2.26%  '

$ws.Range("D3").Value = '3.476.42'
$ws.Range("E3").Value = '  +1.33%  '

$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").Value = '595.10'
$ws.Range("E5").Value = '  +2.06%  '

$ws.Range("D6").Value = '182.34'
$ws.Range("E6").Value = '  +6.26%  '

$ws.Range("D7").Value = '0.616'
$ws.Range("E7").Value = '  +6.57%  '

$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("D9").Value = '3.475.80'
$ws.Range("E9").Value = '  +0.81%  '

$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  +12.09%  '

$ws.Range("D11").Value = '7.00'
$ws.Range("E11").Value = '  -0.75%  '

$ws.Range("D12").Value = '0.431'
$ws.Range("E12").Value = '  +2.97%  '

$ws.Range("D13").Value = '4.071.82'
$ws.Range("E13").Value = '  +0.67%  '

$ws.Range("D14").Value = '32.16'
$ws.Range("E14").Value = '  +7.25%  '

$ws.Range("E15").Value = '  -0.23%  '

$ws.Range("D16").Value = '67.498.92'
$ws.Range("E16").Value = '  +1.90%  '

$ws.Range("D17").Value = '0.0000178'
$ws.Range("E17").Value = '  +3.57%  '

$ws.Range("D18").Value = '3.471.30'
$ws.Range("E18").Value = '  +0.42%  '

$ws.Range("E19").Value = '  +1.13%  '

$ws.Range("D20").Value = '14.18'
$ws.Range("E20").Value = '  +0.66%  '

$ws.Range("D21").Value = '395.15'
$ws.Range("E21").Value = '  +3.67%  '

$ws.Range("D22").Value = '7.97'
$ws.Range("E22").Value = '  +3.69%  '

$ws.Range("D23").Value = '5.80'
$ws.Range("E23").Value = '  +1.71%  '

$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").Value = '0.540'
$ws.Range("E24").Value = '  +3.05%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '0.997'
$ws.Range("E25").Value = '  +0.07%  '

$ws.Range("D26").Value = '71.94'
$ws.Range("E26").Value = '  -0.40%  '

$ws.Range("E27").Value = '  +4.66%  '

$ws.Range("D28").Value = '10.39'
$ws.Range("E28").Value = '  +3.74%  '

$ws.Range("D29").Value = '0.176'
$ws.Range("E29").Value = '  +0.23%  '

$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("D31").Value = '6.15'
$ws.Range("E31").Value = '  +3.07%  '

$ws.Range("E32").Value = '  +1.98%  '

$ws.Range("E33").Value = '  +2.14%  '

$ws.Range("D34").Value = '23.56'
$ws.Range("E34").Value = '  +2.49%  '

$ws.Range("D35").Value = '7.35'
$ws.Range("E35").Value = '  +2.73%  '

$ws.Range("E37").Value = '  +0.25%  '

$ws.Range("D38").Value = '160.97'
$ws.Range("E38").Value = '  -0.79%  '

$ws.Range("D39").Value = '0.891'
$ws.Range("E39").Value = '  +2.96%  '

$ws.Range("D40").Value = '2.87'
$ws.Range("E40").Value = '  +17.45%  '

$ws.Range("E41").Value = '  +0.26%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '6.79'
$ws.Range("E42").Value = '  +0.32%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '4.67'
$ws.Range("E43").Value = '  +3.23%  '

$ws.Range("D44").Value = '26.29'
$ws.Range("E44").Value = '  +2.53%  '

$ws.Range("D45").Value = '0.0721'
$ws.Range("E45").Value = '  +1.66%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '26.44'
$ws.Range("E46").Value = '  -0.03%  '

$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.753.62'
$ws.Range("E47").Value = '  +0.29%  '

$ws.Range("D48").Value = '41.53'
$ws.Range("E48").Value = '  -1.20%  '

$ws.Range("D49").Value = '0.0299'
$ws.Range("E49").Value = '  +2.46%  '

$ws.Range("D50").Value = '326.99'
$ws.Range("E50").Value = '  -0.59%  '

$ws.Range("E51").Value = '  +0.03%  '
